# "Generate Report for Handback" -- refresh the localization-status report:
#  - the handback status text moves from "Ready for handoff" to
#    "Handed back: in sync with en-US" (Overview + both language sheets)
#  - the per-language "Latest Handback DateTime" timestamps advance
#  - the (now resolved) handback-version-mismatch error message is cleared
#    out of the "Error Detail" column for both language sheets

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---- Overview sheet ------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Columns.Item(5).ColumnWidth = 29.17
$overview.Columns.Item(6).ColumnWidth = 29.17

# ---- zh-cn sheet -----------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $statusText
$zhcn.Columns.Item(3).ColumnWidth = 29.17
$zhcn.Range("K2").Value = "2016-09-03 12:53:30"
$zhcn.Range("P2").Value = ""
$zhcn.Columns.Item(16).ColumnWidth = 12.83

# ---- de-de sheet -----------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $statusText
$dede.Columns.Item(3).ColumnWidth = 29.17
$dede.Range("K2").Value = "2016-09-03 12:53:37"
$dede.Range("P2").Value = ""
$dede.Columns.Item(16).ColumnWidth = 12.83
